# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
# (GitHub Actions data-update commit).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $text) {
    # Force the written value to stay a text cell (matches the source data,
    # which stores these as inline strings, not numbers) by using Excel's
    # leading-apostrophe text-entry marker, then reset the style back to
    # Normal so no numeric/text cell format is left behind on the cell.
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "26.245.38"
Set-TextCell $ws.Range("E2") "  +1.74%  "
Set-TextCell $ws.Range("D3") "1.645.79"
Set-TextCell $ws.Range("E3") "  +0.61%  "
Set-TextCell $ws.Range("E4") "  -0.09%  "
Set-TextCell $ws.Range("D5") "217.01"
Set-TextCell $ws.Range("E5") "  +0.71%  "
Set-TextCell $ws.Range("D6") "0.505"
Set-TextCell $ws.Range("E6") "  +0.12%  "
Set-TextCell $ws.Range("E7") "  -0.12%  "
Set-TextCell $ws.Range("E8") "  +0.00%  "
Set-TextCell $ws.Range("D9") "0.0638"
Set-TextCell $ws.Range("E9") "  -0.13%  "
Set-TextCell $ws.Range("D10") "19.99"
Set-TextCell $ws.Range("E10") "  +1.37%  "
Set-TextCell $ws.Range("E11") "  +0.17%  "
Set-TextCell $ws.Range("E12") "  +0.42%  "
Set-TextCell $ws.Range("D13") "1.873.30"
Set-TextCell $ws.Range("E13") "  +0.63%  "
Set-TextCell $ws.Range("D14") "1.642.65"
Set-TextCell $ws.Range("E14") "  +0.43%  "
Set-TextCell $ws.Range("D15") "0.549"
Set-TextCell $ws.Range("E15") "  -2.30%  "
Set-TextCell $ws.Range("E16") "  -0.57%  "
Set-TextCell $ws.Range("E17") "  +0.43%  "
Set-TextCell $ws.Range("D18") "26.231.38"
Set-TextCell $ws.Range("E18") "  +1.64%  "
Set-TextCell $ws.Range("E19") "  -0.13%  "
Set-TextCell $ws.Range("D20") "195.31"
Set-TextCell $ws.Range("E20") "  +1.31%  "
Set-TextCell $ws.Range("D21") "4.43"
Set-TextCell $ws.Range("E21") "  -0.83%  "
Set-TextCell $ws.Range("D23") "6.34"
Set-TextCell $ws.Range("E23") "  -0.58%  "
Set-TextCell $ws.Range("D24") "143.40"
Set-TextCell $ws.Range("E24") "  +0.97%  "
Set-TextCell $ws.Range("E25") "  -0.09%  "
Set-TextCell $ws.Range("E26") "  -1.99%  "
Set-TextCell $ws.Range("E27") "  +1.73%  "
Set-TextCell $ws.Range("E28") "  -0.29%  "
Set-TextCell $ws.Range("D29") "15.61"
Set-TextCell $ws.Range("E29") "  +0.43%  "
Set-TextCell $ws.Range("E30") "  +1.13%  "
Set-TextCell $ws.Range("D31") "0.0504"
Set-TextCell $ws.Range("E31") "  +1.98%  "
Set-TextCell $ws.Range("E32") "  +0.00%  "
Set-TextCell $ws.Range("E34") "  +1.47%  "
Set-TextCell $ws.Range("E35") "  +0.86%  "
Set-TextCell $ws.Range("D36") "0.912"
Set-TextCell $ws.Range("E36") "  +0.71%  "
Set-TextCell $ws.Range("D37") "1.131.45"
Set-TextCell $ws.Range("E37") "  -0.18%  "
Set-TextCell $ws.Range("D38") "0.552"
Set-TextCell $ws.Range("E38") "  +1.26%  "
Set-TextCell $ws.Range("E40") "  +0.91%  "
Set-TextCell $ws.Range("E41") "  -0.13%  "
Set-TextCell $ws.Range("D42") "5.65"
Set-TextCell $ws.Range("E42") "  +1.68%  "
Set-TextCell $ws.Range("E43") "  -0.46%  "
Set-TextCell $ws.Range("E44") "  -1.22%  "
Set-TextCell $ws.Range("D45") "1.782.23"
Set-TextCell $ws.Range("E45") "  +0.64%  "
Set-TextCell $ws.Range("E46") "  +1.56%  "
Set-TextCell $ws.Range("B47") "Cronos"
Set-TextCell $ws.Range("C47") "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell $ws.Range("D47") "0.0517"
Set-TextCell $ws.Range("E47") "  +2.77%  "
Set-TextCell $ws.Range("B48") "RenderToken"
Set-TextCell $ws.Range("C48") "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell $ws.Range("D48") "1.47"
Set-TextCell $ws.Range("E48") "  +3.91%  "
Set-TextCell $ws.Range("B49") "Mantle"
Set-TextCell $ws.Range("C49") "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextCell $ws.Range("D49") "0.418"
Set-TextCell $ws.Range("E49") "  +0.23%  "
Set-TextCell $ws.Range("B50") "EnergySwap"
Set-TextCell $ws.Range("C50") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell $ws.Range("D50") "7.68"
Set-TextCell $ws.Range("E50") "  +2.34%  "
Set-TextCell $ws.Range("E51") "  +1.96%  "
